$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2015, 0, 517758.6681979147, 7332.199167817354, 244.4447294782518, 668.44367586591, 1270),
    @(2015, 0, 517758.6681979147, 7332.199167817354, 244.4447294782518, 668.44367586591, 1270),
    @(2015, 0, 110380.789438057, 1563.148976863473, 52.11308640664471, 142.5052735356513, 240),
    @(2015, 0, 125680.6940147651, 1779.81738726883, 59.33649233935003, 162.2579596494127, 320),
    @(2015, 0, 145523.4041453484, 2060.818385676509, 68.70465207778945, 187.8755589550305, 360),
    @(2015, 0, 145606.7113463998, 2061.998134140415, 68.74398315537958, 187.9831113246758, 370),
    @(2015, 0, 145606.7113463998, 2061.998134140415, 68.74398315537958, 187.9831113246758, 370),
    @(2015, 0, 145606.7113463998, 2061.998134140415, 68.74398315537958, 187.9831113246758, 370),
    @(2015, 0, 173399.4401752912, 2455.582704917485, 81.8655135077677, 223.8644494110184, 520),
    @(2015, 0, 173470.1577400294, 2456.584166216959, 81.89890075363716, 223.9557481411622, 530),
    @(2015, 0, 173470.1577400294, 2456.584166216959, 81.89890075363716, 223.9557481411622, 530),
    @(2015, 0, 223214.9683322945, 3161.041437988354, 105.3844694467967, 288.1779545278223, 720),
    @(2015, 0, 521779.0801168113, 7389.133918962274, 246.3428502907218, 673.6341614852198, 1290),
    @(2015, 0, 529765.3737652395, 7502.231196202831, 250.113347068956, 683.9447324342261, 1290),
    @(2015, 0, 550859.0209042045, 7800.947241913791, 260.0721004133788, 711.1773330589144, 1360),
    @(2015, 0, 550859.0209042045, 7800.947241913791, 260.0721004133788, 711.1773330589144, 1360),
    @(2015, 0, 442612.5120857912, 6268.022714272991, 208.9666526625076, 571.427486849563, 1140),
    @(2015, 0, 64059.51804790787, 907.1738896439488, 30.24384239552926, 82.70297022190482, 240),
    @(2015, 0, 64059.51804790788, 907.1738896439488, 30.24384239552926, 82.70297022190483, 240),
    @(2015, 0, 237529.6757879372, 3363.758055860908, 112.1427431491664, 306.6587183630671, 770),
    @(2015, 0, 410999.8335279666, 5820.342222077867, 194.0416439028035, 530.6144665042294, 1300),
    @(2015, 0, 172522.4580829261, 2451.859218415024, 81.79114449542097, 222.6579828375058, 530),
    @(2015, 0, 172522.4580829261, 2451.859218415024, 81.79114449542097, 222.6579828375058, 530)
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
